# Update the "Förändrad" (Changed) date column (C) for rows 2-27
# from serial date 45205 (2023-10-06) to 45206 (2023-10-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value = 45206
    }
}
